$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'61.405.34"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -1.34%  '

$ws.Range('D3').Value = "'2.985.63"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.28%  '

$ws.Range('D4').Value = "'1.00"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.13%  '

$ws.Range('D5').Value = "'601.34"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.59%  '

$ws.Range('D6').Value = "'143.70"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.69%  '

$ws.Range('E7').Value = '  -0.02%  '

$ws.Range('E8').Value = '  -0.56%  '

$ws.Range('D9').Value = "'2.983.28"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.19%  '

$ws.Range('D10').Value = "'6.06"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +7.71%  '

$ws.Range('E11').Value = '  -1.22%  '

$ws.Range('D12').Value = "'0.455"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +3.04%  '

$ws.Range('E13').Value = '  +0.26%  '

$ws.Range('E14').Value = '  -0.45%  '

$ws.Range('D15').Value = "'0.125"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.38%  '

$ws.Range('D16').Value = "'3.479.30"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.45%  '

$ws.Range('D17').Value = "'6.92"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -2.28%  '

$ws.Range('D18').Value = "'61.395.42"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.47%  '

$ws.Range('D19').Value = "'2.983.86"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.51%  '

$ws.Range('D20').Value = "'449.48"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.86%  '

$ws.Range('D21').Value = "'14.17"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.45%  '

$ws.Range('D22').Value = "'0.684"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.01%  '

$ws.Range('D23').Value = "'7.32"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.60%  '

$ws.Range('D24').Value = "'81.81"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.23%  '

$ws.Range('D25').Value = "'2.20"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.94%  '

$ws.Range('E26').Value = '  +5.20%  '

$ws.Range('E27').Value = '  -2.44%  '

$ws.Range('E28').Value = '  +0.16%  '

$ws.Range('E29').Value = '  +3.20%  '

$ws.Range('E30').Value = '  -0.06%  '

$ws.Range('D31').Value = "'7.13"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.26%  '

$ws.Range('E32').Value = '  -1.55%  '

$ws.Range('D33').Value = "'27.21"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.49%  '

$ws.Range('D34').Value = "'0.109"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.87%  '

$ws.Range('E35').Value = '  +4.71%  '

$ws.Range('D36').Value = "'1.02"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.40%  '

$ws.Range('E37').Value = '  +0.69%  '

$ws.Range('D38').Value = "'50.35"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.63%  '

$ws.Range('E39').Value = '  -2.22%  '

$ws.Range('D40').Value = "'9.09"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.76%  '

$ws.Range('B41').Value = 'dogwifhat'
$ws.Range('C41').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D41').Value = "'2.88"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.10%  '

$ws.Range('B42').Value = 'Kaspa'
$ws.Range('C42').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D42').Value = "'0.122"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +10.25%  '

$ws.Range('D43').Value = "'398.90"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.18%  '

$ws.Range('D44').Value = "'39.59"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +4.18%  '

$ws.Range('E45').Value = '  +0.00%  '

$ws.Range('D46').Value = "'0.269"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.26%  '

$ws.Range('D47').Value = "'2.690.48"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.83%  '

$ws.Range('D48').Value = "'131.48"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +2.87%  '

$ws.Range('E49').Value = '  +0.11%  '

$ws.Range('E50').Value = '  -0.51%  '

$ws.Range('E51').Value = '  +1.09%  '
